$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 542-543, pushing the existing data (rows 542-577)
# down to 544-579.
$ws.Range("A542:R543").Insert()

# Row 542: new "1a (guarda)" Región de O'Higgins entry
$ws.Cells.Item(542, 1).Value = 11
$ws.Cells.Item(542, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(542, 3).Value = "Bíobío"
$ws.Cells.Item(542, 4).Value = 44826
$ws.Cells.Item(542, 5).Value = 8
$ws.Cells.Item(542, 6).Value = 100112004
$ws.Cells.Item(542, 7).Value = "Cebolla"
$ws.Cells.Item(542, 8).Value = "Sin especificar"
$ws.Cells.Item(542, 9).Value = "1a (guarda)"
$ws.Cells.Item(542, 10).Value = 800
$ws.Cells.Item(542, 11).Value = 10000
$ws.Cells.Item(542, 12).Value = 10500
$ws.Cells.Item(542, 13).Value = 10250
$ws.Cells.Item(542, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(542, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(542, 16).Value = 569
$ws.Cells.Item(542, 17).Value = 18
$ws.Cells.Item(542, 18).Value = "Hortaliza"

# Row 543: new "2a (guarda)" Región de O'Higgins entry
$ws.Cells.Item(543, 1).Value = 11
$ws.Cells.Item(543, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(543, 3).Value = "Bíobío"
$ws.Cells.Item(543, 4).Value = 44826
$ws.Cells.Item(543, 5).Value = 8
$ws.Cells.Item(543, 6).Value = 100112004
$ws.Cells.Item(543, 7).Value = "Cebolla"
$ws.Cells.Item(543, 8).Value = "Sin especificar"
$ws.Cells.Item(543, 9).Value = "2a (guarda)"
$ws.Cells.Item(543, 10).Value = 400
$ws.Cells.Item(543, 11).Value = 9500
$ws.Cells.Item(543, 12).Value = 9500
$ws.Cells.Item(543, 13).Value = 9500
$ws.Cells.Item(543, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(543, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(543, 16).Value = 528
$ws.Cells.Item(543, 17).Value = 18
$ws.Cells.Item(543, 18).Value = "Hortaliza"
